# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览": F2, F6, F9, F10, F13, F14, F15, F17, F18
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6416
$ws1.Range("F6").Value  = 57
$ws1.Range("F9").Value  = 84
$ws1.Range("F10").Value = 73
$ws1.Range("F13").Value = 368
$ws1.Range("F14").Value = 786
$ws1.Range("F15").Value = 3131
$ws1.Range("F17").Value = 189
$ws1.Range("F18").Value = 1801

# Sheet "全部类型": F2, F6, F10, F11, F14, F15, F16, F18, F19
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6416
$ws4.Range("F6").Value  = 57
$ws4.Range("F10").Value = 84
$ws4.Range("F11").Value = 73
$ws4.Range("F14").Value = 368
$ws4.Range("F15").Value = 786
$ws4.Range("F16").Value = 3131
$ws4.Range("F18").Value = 189
$ws4.Range("F19").Value = 1801
